$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume(1h)) hold numeric-looking text that must
# stay as plain text, so force text number format before assigning.
$textCells = @("D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6", "D7", "E7", "E8", "D9", "E9", "D10", "E10", "D11", "E11", "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16", "D17", "E17", "D18", "E18", "E19", "E21", "D22", "E22", "D23", "E23", "D24", "E24", "D25", "E25", "E26", "E27", "D40", "E40", "D41", "E41", "D42", "E42", "D43", "E43", "D44", "E44", "D45", "E45", "D46", "E46", "D47", "E47")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated values
$ws.Range("D2").Value = "275.04"
$ws.Range("E2").Value = "-2.28%"
$ws.Range("D3").Value = "27.16"
$ws.Range("E3").Value = "1.40%"
$ws.Range("D4").Value = "4.758"
$ws.Range("E4").Value = "-3.69%"
$ws.Range("D5").Value = "0.06296"
$ws.Range("E5").Value = "-1.79%"
$ws.Range("D6").Value = "6.931"
$ws.Range("E6").Value = "-0.91%"
$ws.Range("D7").Value = "1.345"
$ws.Range("E7").Value = "31.95%"
$ws.Range("E8").Value = "-1.00%"
$ws.Range("D9").Value = "0.1511"
$ws.Range("E9").Value = "1.35%"
$ws.Range("D10").Value = "0.05065"
$ws.Range("E10").Value = "-1.86%"
$ws.Range("D11").Value = "0.07584"
$ws.Range("E11").Value = "2.53%"
$ws.Range("D12").Value = "0.02920"
$ws.Range("E12").Value = "-6.15%"
$ws.Range("D13").Value = "0.08998"
$ws.Range("E13").Value = "-0.61%"
$ws.Range("D14").Value = "0.001561"
$ws.Range("E14").Value = "-1.51%"
$ws.Range("D15").Value = "0.0006348"
$ws.Range("E15").Value = "0.56%"
$ws.Range("D16").Value = "0.005944"
$ws.Range("E16").Value = "-1.45%"
$ws.Range("D17").Value = "3.447"
$ws.Range("E17").Value = "-1.73%"
$ws.Range("D18").Value = "3.303"
$ws.Range("E18").Value = "-1.48%"
$ws.Range("E19").Value = "-1.17%"
$ws.Range("E21").Value = "-0.40%"
$ws.Range("D22").Value = "3.909"
$ws.Range("E22").Value = "-0.83%"
$ws.Range("D23").Value = "0.04396"
$ws.Range("E23").Value = "1.06%"
$ws.Range("D24").Value = "0.001173"
$ws.Range("E24").Value = "-0.30%"
$ws.Range("D25").Value = "0.003836"
$ws.Range("E25").Value = "4.00%"
$ws.Range("E26").Value = "0.12%"
$ws.Range("E27").Value = "14.45%"
$ws.Range("D40").Value = "0.04096"
$ws.Range("E40").Value = "-0.05%"
$ws.Range("D41").Value = "0.006858"
$ws.Range("E41").Value = "3.10%"
$ws.Range("D42").Value = "0.1170"
$ws.Range("E42").Value = "-0.72%"
$ws.Range("D43").Value = "0.002121"
$ws.Range("E43").Value = "-10.07%"
$ws.Range("D44").Value = "0.01152"
$ws.Range("E44").Value = "-11.83%"
$ws.Range("D45").Value = "0.00005179"
$ws.Range("E45").Value = "-1.29%"
$ws.Range("D46").Value = "1.490"
$ws.Range("E46").Value = "-36.76%"
$ws.Range("D47").Value = "0.02302"
$ws.Range("E47").Value = "2.36%"

# Plain text fields (coin name / link) for the BOLO / CoinbaseStockToken swap
$ws.Range("B46").Value = "BOLO"
$ws.Range("C46").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("B47").Value = "CoinbaseStockToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
